# Softball Victoria - Umpiring Development Pathway
# 1) Rename "Development Umpire Clinic" -> "Level 2 Development Umpire Clinic"
#    inside the "Available Training" box (Rectangle 20).
# 2) Re-cache the "datetimeFigureOut" date placeholders (21/10/2013 -> 22/10/2013)
#    on the slide master and every slide layout (mirrors PowerPoint's own
#    re-cache-on-save behaviour for the auto date field).

$p = $ppt.ActivePresentation

# --- 1) Clinic name rename ------------------------------------------------

# The phrase "Development Umpire Clinic" also occurs as a *substring* of the
# unrelated paragraph "Attend Level 2 Development Umpire Clinic" elsewhere on
# the slide, so a plain substring search is not safe. Split each shape's text
# on the paragraph separator (CR, chr(13)) and only touch the paragraph whose
# text is *exactly* "Development Umpire Clinic".

foreach ($slideIdx in 1..$p.Slides.Count) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shapeIdx in 1..$slide.Shapes.Count) {
        $shp = $slide.Shapes.Item($shapeIdx)
        if (-not $shp.HasTextFrame) { continue }
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        $target = "Development Umpire Clinic"
        $parts = $full.Split([char]13)
        $offset = 0
        foreach ($part in $parts) {
            if ($part -eq $target) {
                # Only replace the leading "Development " (12 chars incl.
                # trailing space) with "Level 2 Development " so the
                # paragraph reads "Level 2 Development Umpire Clinic"
                # afterwards, splitting the text into two runs just like the
                # authored edit.
                $sub = $tr.Characters($offset + 1, 12)
                $sub.Text = "Level 2 Development "
            }
            $offset += $part.Length + 1
        }
    }
}

# --- 2) Date placeholder re-cache -----------------------------------------

function Update-DatePlaceholders($shapes) {
    foreach ($i in 1..$shapes.Count) {
        $shp = $shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "21/10/2013") {
            $tr.Text = "22/10/2013"
        }
    }
}

# Slide master
Update-DatePlaceholders $p.SlideMaster.Shapes

# Every slide layout that hangs off the (single) design/master
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
foreach ($layoutIdx in 1..$master.CustomLayouts.Count) {
    $layout = $master.CustomLayouts.Item($layoutIdx)
    Update-DatePlaceholders $layout.Shapes
}

# Notes master (best effort - some runtimes expose this read-only)
try {
    Update-DatePlaceholders $p.NotesMaster.Shapes
} catch {
    Write-Host "NotesMaster date placeholder update skipped:" $_
}
